$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 16-31 (no longer present in the updated table)
$ws.Rows("16:31").Delete()

# New data for rows 2-15 (columns B, C, D, E, F)
$data = @(
    @("", "NSE:ALMONDZ", "NSE:ALKEM", "", ""),
    @("", "NSE:APOLLOPIPE", "", "", ""),
    @("", "NSE:ASHAPURMIN", "", "", ""),
    @("", "NSE:BERGEPAINT", "", "", ""),
    @("", "NSE:BSOFT", "", "", ""),
    @("", "NSE:FIVESTAR", "", "", ""),
    @("", "NSE:GENESYS", "", "", ""),
    @("", "NSE:GFLLIMITED", "", "", ""),
    @("", "NSE:IVZINGOLD", "", "", ""),
    @("", "NSE:KALYANKJIL", "", "", ""),
    @("", "NSE:KSB", "", "", ""),
    @("", "NSE:LICHSGFIN", "", "", ""),
    @("", "NSE:MATRIMONY", "", "", ""),
    @("", "NSE:NDLVENTURE", "", "", "")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
    $ws.Cells.Item($row, 6).Value = $data[$i][4]
}
